$d = $word.ActiveDocument

# Insert "a " before "significant risk" in the first paragraph.
$d.Content.Find.Execute("there is significant risk", $true, $false, $false, $false, $false, $true, 1, $false, "there is a significant risk", 2)

# Fix the comma before "but" in the first paragraph.
$d.Content.Find.Execute("vulnerability scanning, but", $true, $false, $false, $false, $false, $true, 1, $false, "vulnerability scanning but", 2)

# Replace both occurrences of "CISA team" with "security team".
$d.Content.Find.Execute("CISA team", $true, $false, $false, $false, $false, $true, 1, $false, "security team", 2)
